$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update week totals in column B (week 8 results added)
$ws.Range("B2").Value = 2647.6
$ws.Range("B3").Value = 2511.1999999999998
$ws.Range("B4").Value = 2431.6999999999998
$ws.Range("B5").Value = 1870.4
$ws.Range("B6").Value = 1813.9
$ws.Range("B7").Value = 1494.1
$ws.Range("B8").Value = 1240.9000000000001
$ws.Range("B9").Value = 1013.7
$ws.Range("B10").Value = 862.3
$ws.Range("B11").Value = 182

# Add new formatted (but empty) cells in column F, matching the style of column D
$ws.Range("F2:F11").NumberFormat = $ws.Range("D2:D11").NumberFormat

# Update the selection to match the new active cell / selection range
$ws.Range("E2:F11").Select()
